$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 98.56521600000001
$ws.Range("J9").Value = 227.5
$ws.Range("L9").Value = 227.5
$ws.Range("N9").Value = -565.5
$ws.Range("H33").Value = 1055.2941
$ws.Range("I33").Value = 1260
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 1260
$ws.Range("L33").Value = 100
$ws.Range("M33").Value = -1031
$ws.Range("N33").Value = -558
$ws.Range("H113").Value = 7145643
$ws.Range("I113").Value = 10002000
$ws.Range("J113").Value = 4749.75
$ws.Range("K113").Value = 10002000
$ws.Range("L113").Value = 4749.75
$ws.Range("M113").Value = -9998746
$ws.Range("N113").Value = -11257.75
$ws.Range("H141").Value = 2793.125
$ws.Range("I141").Value = 695.125
$ws.Range("J141").Value = 9087.125
$ws.Range("K141").Value = 2085.375
$ws.Range("L141").Value = 27261.375
$ws.Range("M141").Value = 3094.625
$ws.Range("N141").Value = -37621.375
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3256.03
$ws.Range("I32").Value = 3256.03
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3256.03
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -2969.03
$ws.Range("H37").Value = 9080.143
$ws.Range("I37").Value = 3836.7778
$ws.Range("J37").Value = 13012.667
$ws.Range("K37").Value = 3836.7778
$ws.Range("L37").Value = 13012.667
$ws.Range("M37").Value = -3563.7778
$ws.Range("N37").Value = -13558.667
$ws.Range("H61").Value = 2957045.5
$ws.Range("I61").Value = 1544448.1
$ws.Range("J61").Value = 8405635
$ws.Range("K61").Value = 1544448.1
$ws.Range("L61").Value = 8405635
$ws.Range("M61").Value = -1544236.1
$ws.Range("N61").Value = -8406059
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").ClearContents()
$ws.Range("N76").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").ClearContents()
$ws.Range("N79").Value = 0
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = 0
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = 0
$ws.Range("H136").Value = 2957045.5
$ws.Range("I136").Value = 1544448.1
$ws.Range("J136").Value = 8405635
$ws.Range("K136").Value = 4633344.300000001
$ws.Range("L136").Value = 25216905
$ws.Range("M136").Value = -4630794.300000001
$ws.Range("N136").Value = -25222005
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 8174751
$ws.Range("I64").Value = 454685.62
$ws.Range("J64").Value = 16666823
$ws.Range("K64").Value = 454685.62
$ws.Range("L64").Value = 16666823
$ws.Range("M64").Value = -454460.62
$ws.Range("N64").Value = -16667273
$ws.Range("H67").Value = 8174751
$ws.Range("I67").Value = 454685.62
$ws.Range("J67").Value = 16666823
$ws.Range("K67").Value = 454685.62
$ws.Range("L67").Value = 16666823
$ws.Range("M67").Value = -453905.62
$ws.Range("N67").Value = -16668383
$ws.Range("H76").Value = 32000
$ws.Range("J76").Value = 32000
$ws.Range("L76").Value = 32000
$ws.Range("N76").Value = -32630
$ws.Range("H79").Value = 32000
$ws.Range("J79").Value = 32000
$ws.Range("L79").Value = 32000
$ws.Range("N79").Value = -34184
$ws.Range("H107").Value = 1001084.9
$ws.Range("I107").Value = 1001084.9
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1001084.9
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -999164.9
$ws.Range("H134").Value = 12079905
$ws.Range("I134").Value = 13655708
$ws.Range("J134").Value = 103800
$ws.Range("K134").Value = 40967124
$ws.Range("L134").Value = 311400
$ws.Range("M134").Value = -40964589
$ws.Range("N134").Value = -316470
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9657.111000000001
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 9657.111000000001
$ws.Range("K31").Value = 0
$ws.Range("L31").ClearContents()
$ws.Range("M31").Value = 9657.111000000001
$ws.Range("N31").Value = -10247.111
$ws.Range("H34").Value = 9657.111000000001
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 9657.111000000001
$ws.Range("K34").Value = 0
$ws.Range("L34").ClearContents()
$ws.Range("M34").Value = 9657.111000000001
$ws.Range("N34").Value = -10061.111
$ws.Range("H50").Value = 12221.6
$ws.Range("J50").Value = 12221.6
$ws.Range("L50").Value = 12221.6
$ws.Range("N50").Value = -13471.6
$ws.Range("H51").Value = 28030.1
$ws.Range("J51").Value = 10037.625
$ws.Range("L51").Value = 10037.625
$ws.Range("N51").Value = -11509.625
$ws.Range("H59").Value = 17396.2
$ws.Range("J59").Value = 17396.2
$ws.Range("L59").Value = 17396.2
$ws.Range("N59").Value = -19686.2
$ws.Range("H60").Value = 15742.267
$ws.Range("J60").Value = 9613.4
$ws.Range("L60").Value = 9613.4
$ws.Range("N60").Value = -10635.4
$ws.Range("H61").Value = 28030.1
$ws.Range("J61").Value = 10037.625
$ws.Range("L61").Value = 10037.625
$ws.Range("N61").Value = -10733.625
$ws.Range("H68").Value = 17888
$ws.Range("J68").Value = 17888
$ws.Range("L68").Value = 17888
$ws.Range("N68").Value = -19386
$ws.Range("H71").Value = 17888
$ws.Range("J71").Value = 17888
$ws.Range("L71").Value = 53664
$ws.Range("N71").Value = -61152
$ws.Range("H74").Value = 16446.076
$ws.Range("J74").Value = 17709.5
$ws.Range("L74").Value = 17709.5
$ws.Range("N74").Value = -19457.5
$ws.Range("H77").Value = 16446.076
$ws.Range("J77").Value = 17709.5
$ws.Range("L77").Value = 53128.5
$ws.Range("N77").Value = -61864.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3099604
$ws.Range("I5").Value = 4525378
$ws.Range("J5").Value = 2233955.8
$ws.Range("K5").Value = 13576134
$ws.Range("L5").Value = 6701867.399999999
$ws.Range("M5").Value = -13576022
$ws.Range("N5").Value = -6702091.399999999
$ws.Range("H112").Value = 6099
$ws.Range("I112").Value = 2475.6667
$ws.Range("J112").Value = 7186
$ws.Range("K112").Value = 7427.000100000001
$ws.Range("L112").Value = 21558
$ws.Range("M112").Value = -6319.000100000001
$ws.Range("N112").Value = -23774
$ws.Range("H135").Value = 3099604
$ws.Range("I135").Value = 4525378
$ws.Range("J135").Value = 2233955.8
$ws.Range("K135").Value = 40728402
$ws.Range("L135").Value = 20105602.2
$ws.Range("M135").Value = -40725867
$ws.Range("N135").Value = -20110672.2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 23792.916
$ws.Range("I113").Value = 1057.2222
$ws.Range("J113").Value = 92000
$ws.Range("K113").Value = 1057.2222
$ws.Range("L113").Value = 92000
$ws.Range("M113").Value = 1112.7778
$ws.Range("N113").Value = -96340
$ws.Range("H122").Value = 3692.2
$ws.Range("I122").Value = 2796.7896
$ws.Range("J122").Value = 5238.8184
$ws.Range("K122").Value = 8390.3688
$ws.Range("L122").Value = 15716.4552
$ws.Range("M122").Value = -5940.3688
$ws.Range("N122").Value = -20616.4552
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 15584.895
$ws.Range("I93").Value = 3756.5557
$ws.Range("J93").Value = 26230.4
$ws.Range("K93").Value = 3756.5557
$ws.Range("L93").Value = 26230.4
$ws.Range("M93").Value = -2508.5557
$ws.Range("N93").Value = -28726.4
$ws.Range("H132").Value = 3764964.8
$ws.Range("I132").Value = 4766468.5
$ws.Range("J132").Value = 9324.875
$ws.Range("K132").Value = 14299405.5
$ws.Range("L132").Value = 27974.625
$ws.Range("M132").Value = -14296875.5
$ws.Range("N132").Value = -33034.625
$ws.Range("H134").Value = 45928.5
$ws.Range("J134").Value = 45928.5
$ws.Range("L134").Value = 45928.5
$ws.Range("N134").Value = -56068.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 29450.625
$ws.Range("J135").Value = 29450.625
$ws.Range("L135").Value = 29450.625
$ws.Range("N135").Value = -39590.625
